$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'94.487.24"
$ws.Range("E2").Value = "  +3.09%  "

$ws.Range("D3").Value = "'3.140.35"
$ws.Range("E3").Value = "  +1.39%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'241.21"
$ws.Range("E5").Value = "  +0.51%  "

$ws.Range("D6").Value = "'619.81"
$ws.Range("E6").Value = "  +0.70%  "

$ws.Range("D7").Value = "'1.14"
$ws.Range("E7").Value = "  +5.18%  "

$ws.Range("D8").Value = "'0.395"
$ws.Range("E8").Value = "  +1.01%  "

$ws.Range("E9").Value = "  -0.06%  "

$ws.Range("D10").Value = "'0.840"
$ws.Range("E10").Value = "  +15.13%  "

$ws.Range("D11").Value = "'3.138.80"
$ws.Range("E11").Value = "  +1.43%  "

$ws.Range("D12").Value = "'0.199"
$ws.Range("E12").Value = "  -1.50%  "

$ws.Range("D13").Value = "'0.0000247"
$ws.Range("E13").Value = "  -1.56%  "

$ws.Range("D14").Value = "'35.46"
$ws.Range("E14").Value = "  +4.19%  "

$ws.Range("D15").Value = "'94.103.27"
$ws.Range("E15").Value = "  +2.33%  "

$ws.Range("D16").Value = "'5.45"
$ws.Range("E16").Value = "  -1.29%  "

$ws.Range("D17").Value = "'3.700.37"
$ws.Range("E17").Value = "  +0.70%  "

$ws.Range("D18").Value = "'3.089.59"
$ws.Range("E18").Value = "  -0.59%  "

$ws.Range("D19").Value = "'3.72"
$ws.Range("E19").Value = "  +2.51%  "

$ws.Range("D20").Value = "'15.27"
$ws.Range("E20").Value = "  +3.80%  "

$ws.Range("D21").Value = "'6.04"
$ws.Range("E21").Value = "  +3.66%  "

$ws.Range("D22").Value = "'452.30"
$ws.Range("E22").Value = "  +2.15%  "

$ws.Range("D23").Value = "'0.0000204"
$ws.Range("E23").Value = "  +1.94%  "

$ws.Range("D24").Value = "'9.16"
$ws.Range("E24").Value = "  -1.51%  "

$ws.Range("E25").Value = "  +4.23%  "

$ws.Range("D26").Value = "'5.75"
$ws.Range("E26").Value = "  +2.49%  "

$ws.Range("D27").Value = "'12.39"
$ws.Range("E27").Value = "  +6.54%  "

$ws.Range("D28").Value = "'86.49"
$ws.Range("E28").Value = "  -0.06%  "

$ws.Range("D29").Value = "'3.278.29"
$ws.Range("E29").Value = "  +0.24%  "

$ws.Range("E30").Value = "  -0.12%  "

$ws.Range("D31").Value = "'0.255"
$ws.Range("E31").Value = "  +12.10%  "

$ws.Range("D32").Value = "'0.182"
$ws.Range("E32").Value = "  +9.34%  "

$ws.Range("E33").Value = "  -7.46%  "

$ws.Range("D34").Value = "'9.50"
$ws.Range("E34").Value = "  +2.57%  "

$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = "  -0.31%  "

$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D36").Value = "'8.10"
$ws.Range("E36").Value = "  +2.44%  "

$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "'0.163"
$ws.Range("E37").Value = "  -3.70%  "

$ws.Range("D38").Value = "'26.55"
$ws.Range("E38").Value = "  +1.55%  "

$ws.Range("D39").Value = "'1.93"
$ws.Range("E39").Value = "  -0.22%  "

$ws.Range("D40").Value = "'0.461"
$ws.Range("E40").Value = "  +7.02%  "

$ws.Range("D41").Value = "'489.76"
$ws.Range("E41").Value = "  +2.47%  "

$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D42").Value = "'1.31"
$ws.Range("E42").Value = "  +0.13%  "

$ws.Range("B43").Value = "MantraDAO"
$ws.Range("C43").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D43").Value = "'3.80"
$ws.Range("E43").Value = "  -6.84%  "

$ws.Range("B44").Value = "WhiteBITCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D44").Value = "'23.97"
$ws.Range("E44").Value = "  +8.09%  "

$ws.Range("D45").Value = "'3.32"
$ws.Range("E45").Value = "  -1.58%  "

$ws.Range("D47").Value = "'159.00"
$ws.Range("E47").Value = "  +0.04%  "

$ws.Range("D48").Value = "'0.707"
$ws.Range("E48").Value = "  +1.60%  "

$ws.Range("D49").Value = "'1.89"
$ws.Range("E49").Value = "  -0.63%  "

$ws.Range("D50").Value = "'4.51"
$ws.Range("E50").Value = "  +3.15%  "

$ws.Range("D51").Value = "'1.34"
$ws.Range("E51").Value = "  -0.97%  "
